$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh headline metrics now that trade #3 has closed.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.94   # Current Capital
$wsSummary.Range("B4").Value = -0.06     # Total P&L $
$wsSummary.Range("B5").Value = -0.4      # Total P&L %
$wsSummary.Range("B6").Value = 3         # Total Trades
$wsSummary.Range("B7").Value = 2         # Winning Trades
$wsSummary.Range("B9").Value = 66.67     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.94      # Capital
$wsStatus.Range("D4").Value = 3          # Trades
$wsStatus.Range("E4").Value = -0.06      # P&L $
$wsStatus.Range("F4").Value = -0.06      # P&L %
$wsStatus.Range("G4").Value = 66.67      # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new trade #3 row (row 4) to a trade log sheet that
# already has the same header/row-2/row-3 layout ("All Trades" & "MarketMaking").
# ---------------------------------------------------------------------------
function Add-TradeRow4($ws) {
    $ws.Range("A4").Value = 3
    # "2026-02-17" looks like a date to Excel's auto-detection, so it must be
    # entered with a leading apostrophe to keep it as literal text, matching
    # the plain string already used in rows 2-3. Style is reset afterwards so
    # no stray quote-prefix formatting is left behind on the cell.
    $ws.Range("B4").Value = "'2026-02-17"
    $ws.Range("B4").Style = "Normal"
    $ws.Range("C4").Value = "07:52:17"
    $ws.Range("D4").Value = "MarketMaking"
    $ws.Range("E4").Value = "DOWN"
    $ws.Range("F4").Value = 0.29
    $ws.Range("G4").Value = 0.33
    $ws.Range("H4").Value = "CLOSED"
    $ws.Range("I4").Value = 13.7931
    $ws.Range("J4").Value = 0.04
    $ws.Range("K4").Value = 99.94
    $ws.Range("L4").Value = 0
    $ws.Range("M4").Value = 0
    $ws.Range("N4").Value = 0.6
    $ws.Range("O4").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P4").Value = "early_exit"
    $ws.Range("Q4").Value = 0.1
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow4 $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow4 $wsMarketMaking
